# Update the NSY (NBA player roster) sheet with the new data set.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New full table (header stays the same, rows 2-18 are the player data).
$data = @(
    @("Shai Gilgeous-Alexander", "PG", "Oklahoma City Thunder"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("T.J. McConnell", "PG", "Indiana Pacers"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("RJ Barrett", "SF,PF", "Toronto Raptors"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Kris Dunn", "PG,SG", "LA Clippers"),
    @("Dennis Schröder", "PG", "Brooklyn Nets"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Jimmy Butler", "SF,PF", "Miami Heat")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
